$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EventRow($Row, $DateSerial, $Event, $Location, $City, $Link) {
    # Date (column A) keeps its existing date-formatted style (s="4").
    $ws.Cells.Item($Row, 1).Value = $DateSerial

    # Plain text columns (B-D) - force the shared "text" style (numFmt "@")
    # that the rest of the table uses, matching the other populated rows.
    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 3).NumberFormat = "@"
    $ws.Cells.Item($Row, 4).NumberFormat = "@"

    $ws.Cells.Item($Row, 2).Value = $Event
    $ws.Cells.Item($Row, 3).Value = $Location
    $ws.Cells.Item($Row, 4).Value = $City

    # Link column (E): add a real hyperlink, then restore the shared
    # text/number format + base font so the cell keeps the plain table
    # style instead of the hyperlink-only style that .Hyperlinks.Add
    # applies automatically.
    $linkCell = $ws.Cells.Item($Row, 5)
    $ws.Hyperlinks.Add($linkCell, $Link, "", "", $Link)

    $linkCell.NumberFormat = "@"
    $linkCell.Font.Name = "Calibri"
    $linkCell.Font.Size = 11
    $linkCell.Font.Underline = $false
    $linkCell.Font.Color = 0

    # Re-apply the underline/blue "link" look directly on the characters so
    # it is baked into the text run (matches how the existing link cells in
    # this workbook are formatted) instead of living only on the cell style.
    $len = $Link.Length
    $run1 = $linkCell.Characters(1, $len - 1)
    $run1.Font.Underline = $true
    $run1.Font.Color = 65280
    $run1.Font.Size = 11
    $run1.Font.Name = "Calibri"

    $run2 = $linkCell.Characters($len, 1)
    $run2.Font.Underline = $true
    $run2.Font.Color = 65280
    $run2.Font.Size = 11
    $run2.Font.Name = "Calibri"
}

Set-EventRow 310 45766 "DAY & NIGHT - OUTDOOR & INDOOR" "SNRS" "Dortmund" "https://www.instagram.com/reel/DIOD77wNrru/?igsh=MXV2ejZkMmRkZ3EzMg=="
Set-EventRow 311 45763 "BEATS BASS COLOGNE EASTER SPECIAL" "Odonien" "Köln" "https://www.instagram.com/p/DHq7HpKswp2/?igsh=MW4waW94eHgyY21mdQ=="
Set-EventRow 312 45773 "TAGESRAVER" "Elektroküche" "Köln" "https://www.instagram.com/reel/DH0sb-hMmus/?igsh=OXc0OTU0MnpwaTVi"
Set-EventRow 313 45759 "DEADLY SINS" "SNRS" "Dortmund" "https://www.instagram.com/reel/DIJsSH2K81q/?igsh=MWVnb21saHpycm4yZg=="
Set-EventRow 314 45773 "NEXORA" "Stollen134" "Dortmund" "https://www.instagram.com/nexora.raves?igsh=dW94ZHlkbzdiYjE4"

Write-Output "rows 310-314 populated"
